# alejandros_loggbok.docx - "forgort to push this before"
#
# Semantic changes in the target diff (ignoring the cosmetic w:proofErr
# spell-check markers that real Word re-inserts on every save and which
# the Word object model has no API to create by hand):
#
#   1) Paragraph "30:de oktober : "                       -> text unchanged
#   2) Paragraph "skapa index.html ... fontawesome."       -> text unchanged
#   3) Paragraph "skapa en footer ... styla den."           -> text unchanged,
#      but the _GoBack bookmark that used to sit in this paragraph is removed
#   4) Paragraph "skappa och styla carousel för huvudsidan" -> text grows with
#      " samt integrera footerns css med index.css och göra båda carousel och
#      footern responsive", and the _GoBack bookmark re-appears at the very
#      end of this (now longer) paragraph, right before the paragraph mark.

$d = $word.ActiveDocument

# --- 1) Move the _GoBack bookmark out of the "footer" paragraph -----------
# (it will be re-created at the end of the "carousel" paragraph below)
$goBack = $null
foreach ($b in $d.Bookmarks) {
    if ($b.Name -eq "_GoBack") { $goBack = $b }
}
if ($goBack -eq $null) {
    # Bookmarks whose name starts with "_" are "hidden" and do not show up
    # while iterating the collection, but can still be addressed by name.
    $goBack = $d.Bookmarks("_GoBack")
}
if ($goBack -ne $null) {
    $goBack.Delete()
}

# --- 2) Grow the last ("carousel") paragraph with the new sentence --------
$lastParaIndex = $d.Paragraphs.Count
$carousel = $d.Paragraphs($lastParaIndex)

# Range covering the paragraph's text, excluding the trailing paragraph mark.
$body = $carousel.Range.Duplicate()
$body.MoveEnd(1, -1)

$newTail = " samt integrera footerns css med index.css och göra båda carousel och footern responsive"
$body.InsertAfter($newTail)

# --- 3) Re-insert the _GoBack bookmark at the very end of that paragraph --
# A zero-length bookmark placed exactly at "end of paragraph text" needs a
# one-character scratch anchor after it, otherwise the anchor collapses to
# the start of the document; so insert a throw-away character, bookmark the
# (now safely interior) boundary in front of it, then remove the character.
$endPos = $body.End
$scratch = $d.Range($endPos, $endPos)
$scratch.InsertAfter("#")

$anchor = $d.Range($endPos, $endPos + 1)
$anchor.Collapse(1)            # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $anchor)

$d.Range($endPos, $endPos + 1).Delete()
